$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 159 (pushes the existing rows 159-230 down
# to 160-231, matching the diff: dimension grows from A1:R230 to A1:R231).
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row with the new weekly price-report entry.
$ws.Cells.Item(159, 1).Value = 10
$ws.Cells.Item(159, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(159, 3).Value = "La Araucanía"
$ws.Cells.Item(159, 4).Value = 45202
$ws.Cells.Item(159, 5).Value = 9
$ws.Cells.Item(159, 6).Value = 100112031
$ws.Cells.Item(159, 7).Value = "Poroto verde"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 65
$ws.Cells.Item(159, 11).Value = 30000
$ws.Cells.Item(159, 12).Value = 30000
$ws.Cells.Item(159, 13).Value = 30000
$ws.Cells.Item(159, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(159, 15).Value = "Perú"
$ws.Cells.Item(159, 16).Value = 1200
$ws.Cells.Item(159, 17).Value = 25
$ws.Cells.Item(159, 18).Value = "Hortaliza"
